$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 9968578
$ws.Range("C3").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 32422
$ws.Range("H3").Value = 32422
